# Updated symbol list edit script
# Applies cell text updates to Sheet1 (cryptos.xlsx) while preserving the
# "text" cell storage (inline/shared string) that the source file uses
# for every data cell, even numeric-looking ones like prices and hours.

function Set-CellText {
    param(
        $ws,
        [string]$addr,
        [string]$text
    )
    $cell = $ws.Range($addr)
    # Force text storage so numeric-looking strings (e.g. "236.89", "17")
    # are not silently re-interpreted as numbers by Excel's smart entry.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    # Drop back to the default style so we don't leave a stray "Text" style
    # applied to cells that originally had no explicit style.
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-CellText $ws 'D2' '236.89'
Set-CellText $ws 'G2' '17'
Set-CellText $ws 'D3' '21.95'
Set-CellText $ws 'G3' '17'
Set-CellText $ws 'D4' '5.360'
Set-CellText $ws 'G4' '17'
Set-CellText $ws 'G5' '17'
Set-CellText $ws 'D6' '6.477'
Set-CellText $ws 'G6' '17'
Set-CellText $ws 'D7' '3.346'
Set-CellText $ws 'G7' '17'
Set-CellText $ws 'D8' '0.7990'
Set-CellText $ws 'G8' '17'
Set-CellText $ws 'D9' '1.039'
Set-CellText $ws 'G9' '17'
Set-CellText $ws 'D10' '0.1384'
Set-CellText $ws 'G10' '17'
Set-CellText $ws 'D11' '0.07310'
Set-CellText $ws 'G11' '17'
Set-CellText $ws 'D12' '0.03159'
Set-CellText $ws 'G12' '17'
Set-CellText $ws 'D13' '0.02950'
Set-CellText $ws 'G13' '17'
Set-CellText $ws 'D14' '0.09246'
Set-CellText $ws 'G14' '17'
Set-CellText $ws 'D15' '0.001667'
Set-CellText $ws 'G15' '17'
Set-CellText $ws 'D16' '3.259'
Set-CellText $ws 'G16' '17'
Set-CellText $ws 'D17' '0.04780'
Set-CellText $ws 'G17' '17'
Set-CellText $ws 'D18' '0.0005716'
Set-CellText $ws 'E18' '17OneONE'
Set-CellText $ws 'G18' '17'
Set-CellText $ws 'D19' '0.006253'
Set-CellText $ws 'G19' '17'
Set-CellText $ws 'D20' '0.005063'
Set-CellText $ws 'G20' '17'
Set-CellText $ws 'D21' '0.001050'
Set-CellText $ws 'G21' '17'
Set-CellText $ws 'D22' '0.0001501'
Set-CellText $ws 'G22' '17'
Set-CellText $ws 'D23' '0.0003703'
Set-CellText $ws 'G23' '17'
Set-CellText $ws 'D24' '3.953'
Set-CellText $ws 'G24' '17'
Set-CellText $ws 'G25' '17'
Set-CellText $ws 'G26' '17'
Set-CellText $ws 'G27' '17'
Set-CellText $ws 'G28' '17'
Set-CellText $ws 'G29' '17'
Set-CellText $ws 'G30' '17'
Set-CellText $ws 'G31' '17'
Set-CellText $ws 'G32' '17'
Set-CellText $ws 'G33' '17'
Set-CellText $ws 'G34' '17'
Set-CellText $ws 'G35' '17'
Set-CellText $ws 'G36' '17'
Set-CellText $ws 'G37' '17'
Set-CellText $ws 'G38' '17'
Set-CellText $ws 'G39' '17'
Set-CellText $ws 'D40' '0.04104'
Set-CellText $ws 'G40' '17'
Set-CellText $ws 'D41' '0.007034'
Set-CellText $ws 'G41' '17'
Set-CellText $ws 'B42' 'BKEXToken'
Set-CellText $ws 'C42' 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
Set-CellText $ws 'D42' '0.1040'
Set-CellText $ws 'E42' '41BKEXTokenBKK'
Set-CellText $ws 'G42' '17'
Set-CellText $ws 'B43' 'CEJI'
Set-CellText $ws 'C43' 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
Set-CellText $ws 'D43' '0.002962'
Set-CellText $ws 'E43' '42CEJICEJI'
Set-CellText $ws 'G43' '17'
Set-CellText $ws 'D44' '0.008807'
Set-CellText $ws 'G44' '17'
Set-CellText $ws 'D45' '0.00005437'
Set-CellText $ws 'G45' '17'
Set-CellText $ws 'D46' '0.00000000751'
Set-CellText $ws 'G46' '17'
Set-CellText $ws 'D47' '0.6759'
Set-CellText $ws 'G47' '17'
Set-CellText $ws 'D48' '0.03617'
Set-CellText $ws 'E48' '47BOLOBOLOWorstin24h'
Set-CellText $ws 'G48' '17'
Set-CellText $ws 'D49' '0.00002102'
Set-CellText $ws 'G49' '17'
Set-CellText $ws 'D50' '0.01011'
Set-CellText $ws 'G50' '17'
Set-CellText $ws 'G51' '17'
